# Apply crypto price/volume updates described by the commit diff.
# Cells in column D/E hold free-form text (prices formatted with
# '.' thousand separators, percentages with padding spaces); force
# Text number-format first so Excel doesn't reinterpret them as numbers
# and silently normalise (e.g. '1.00' -> 1, '9.88' -> 9.88 without the
# trailing-zero formatting, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.623.80'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.22%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.472.10'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.26%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.75'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '92.37'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.28%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.02'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.65%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +7.43%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.851.45'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.53%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.57%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.63'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.484.33'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.71%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.792'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.81%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.591.09'

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.54%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0949'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.35%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.22'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.34'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.83'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.28%  '

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.42%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.97%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.66'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.61%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.27'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.28%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.88'

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.21'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.50%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '161.37'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.32%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.93%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.13%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.40%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.28'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.77%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.74%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.103'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.38%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.98'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.76%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.79%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.987.86'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.85%  '

$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0285'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.99'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.33%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.69%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.22'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.68%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.707.18'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.58%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.57'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.24'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.39%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.21'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.59%  '
